# Move GABRIEL's (account 005666419) balance row from its old position
# further down the sheet (sorted-by-balance list) up near the top, and
# update the balance value from 16.95 to 27016.95.
#
# Data layout: col A = Conta (account, text w/ leading zeros),
#              col B = Nome (text), col C = Saldo (number).
# Row 1 is the header ("Conta"/"Nome"/"Saldo"); the rows below are sorted
# by descending Saldo. GABRIEL currently sits at row 265 with Saldo 16.95;
# the new entry (Saldo 27016.95) belongs right above row 6
# (005064129 / THIAGO / 20357.1), which is its correct sorted position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove GABRIEL's old row (account 005666419, balance 16.95).
$ws.Rows.Item(265).Delete()

# 2) Insert a fresh row just above row 6 and populate the updated entry.
$ws.Rows.Item(6).Insert()

# Leading zeros in the account number must survive as text, not become
# the number 5666419 - prefix with an apostrophe to force text entry,
# then clear the resulting "quote prefix" formatting so the cell keeps
# the plain (unstyled) look of its neighboring account cells.
$ws.Cells.Item(6, 1).Value = "'005666419"
$ws.Cells.Item(6, 1).ClearFormats()

$ws.Cells.Item(6, 2).Value = "GABRIEL"
$ws.Cells.Item(6, 3).Value = 27016.95
